$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3 duplicates row 2's match-summary data (as scraped) for Sheldon Cottrell.
# Columns G:K hold numeric-looking values that must stay stored as text (matching
# the existing numberStoredAsText handling used for row 2), so force a text
# number format before assigning them.
$ws.Range("G3:K3").NumberFormat = "@"

$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " October 08 2020"
$ws.Range("C3").Value = "Sunrisers won by 69 runs"
$ws.Range("D3").Value = "Kings XI Punjab"
$ws.Range("E3").Value = "Sunrisers Hyderabad"
$ws.Range("F3").Value = "Sheldon Cottrell "
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "2"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "0.00"
